{"js": "// Update the thresholds_summary table with refined counts, and update the\n// \"Chosen thresholds\" summary paragraphs to match the new Full/Core-Level\n// analysis selections.\n\nconst body = context.document.body;\n\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// (row, col, newValue) \u2014 row/col are 0-based and include the header row/col.\n// Columns: 0=threshold, 1=row_count, 2=unique_count, 3=repeated_count, 4=no_repeats_bool\nconst cellUpdates = [\n  [10, 1, \"5\"],   // threshold 0.45: row_count 4 -> 5\n  [10, 2, \"5\"],   // threshold 0.45: unique_count 4 -> 5\n  [13, 1, \"7\"],   // threshold 0.6:  row_count 8 -> 7\n  [13, 2, \"7\"],   // threshold 0.6:  unique_count 8 -> 7\n  [14, 1, \"7\"],   // threshold 0.65: row_count 8 -> 7\n  [14, 2, \"7\"],   // threshold 0.65: unique_count 8 -> 7\n  [15, 1, \"7\"],   // threshold 0.7:  row_count 11 -> 7\n  [15, 2, \"7\"],   // threshold 0.7:  unique_count 11 -> 7\n  [16, 1, \"7\"],   // threshold 0.75: row_count 13 -> 7\n  [16, 2, \"7\"],   // threshold 0.75: unique_count 13 -> 7\n  [17, 1, \"7\"],   // threshold 0.8:  row_count 13 -> 7\n  [17, 2, \"7\"],   // threshold 0.8:  unique_count 13 -> 7\n  [18, 1, \"7\"],   // threshold 0.85: row_count 17 -> 7\n  [18, 2, \"7\"],   // threshold 0.85: unique_count 17 -> 7\n  [19, 1, \"7\"],   // threshold 0.9:  row_count 17 -> 7\n  [19, 2, \"7\"],   // threshold 0.9:  unique_count 17 -> 7\n  [20, 1, \"9\"],   // threshold 0.95: row_count 17 -> 9\n  [20, 2, \"7\"],   // threshold 0.95: unique_count 17 -> 7\n  [20, 3, \"4\"],   // threshold 0.95: repeated_count 0 -> 4\n  [20, 4, \"False\"], // threshold 0.95: no_repeats_bool True -> False\n  [21, 1, \"10\"],  // threshold 1.0:  row_count 18 -> 10\n  [21, 2, \"7\"],   // threshold 1.0:  unique_count 17 -> 7\n  [21, 3, \"5\"],   // threshold 1.0:  repeated_count 2 -> 5\n];\n\nfor (const [row, col, value] of cellUpdates) {\n  table.getCell(row, col).value = value;\n}\n\nawait context.sync();\n\n// Update the two \"Chosen thresholds\" paragraphs.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text;\n  if (text === \"  - Full Analysis threshold = 0.95.\") {\n    paragraph.getRange().insertText(\"  - Full Analysis threshold = 0.90.\", \"Replace\");\n  } else if (text === \"  - Core-Level Analysis threshold = 0.65.\") {\n    paragraph.getRange().insertText(\n      \"  - Core-Level Analysis: NONE found (no threshold yields no repeats & \\u22649 rows?).\",\n      \"Replace\"\n    );\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the thresholds_summary table with refined counts, and update the\n# \"Chosen thresholds\" summary paragraphs to match the new Full/Core-Level\n# analysis selections.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Word COM table cells are 1-based (row, col).\n# Columns: 1=threshold, 2=row_count, 3=unique_count, 4=repeated_count, 5=no_repeats_bool\n$t.Cell(11, 2).Range.Text = \"5\"    # threshold 0.45: row_count 4 -> 5\n$t.Cell(11, 3).Range.Text = \"5\"    # threshold 0.45: unique_count 4 -> 5\n\n$t.Cell(14, 2).Range.Text = \"7\"    # threshold 0.6: row_count 8 -> 7\n$t.Cell(14, 3).Range.Text = \"7\"    # threshold 0.6: unique_count 8 -> 7\n\n$t.Cell(15, 2).Range.Text = \"7\"    # threshold 0.65: row_count 8 -> 7\n$t.Cell(15, 3).Range.Text = \"7\"    # threshold 0.65: unique_count 8 -> 7\n\n$t.Cell(16, 2).Range.Text = \"7\"    # threshold 0.7: row_count 11 -> 7\n$t.Cell(16, 3).Range.Text = \"7\"    # threshold 0.7: unique_count 11 -> 7\n\n$t.Cell(17, 2).Range.Text = \"7\"    # threshold 0.75: row_count 13 -> 7\n$t.Cell(17, 3).Range.Text = \"7\"    # threshold 0.75: unique_count 13 -> 7\n\n$t.Cell(18, 2).Range.Text = \"7\"    # threshold 0.8: row_count 13 -> 7\n$t.Cell(18, 3).Range.Text = \"7\"    # threshold 0.8: unique_count 13 -> 7\n\n$t.Cell(19, 2).Range.Text = \"7\"    # threshold 0.85: row_count 17 -> 7\n$t.Cell(19, 3).Range.Text = \"7\"    # threshold 0.85: unique_count 17 -> 7\n\n$t.Cell(20, 2).Range.Text = \"7\"    # threshold 0.9: row_count 17 -> 7\n$t.Cell(20, 3).Range.Text = \"7\"    # threshold 0.9: unique_count 17 -> 7\n\n$t.Cell(21, 2).Range.Text = \"9\"    # threshold 0.95: row_count 17 -> 9\n$t.Cell(21, 3).Range.Text = \"7\"    # threshold 0.95: unique_count 17 -> 7\n$t.Cell(21, 4).Range.Text = \"4\"    # threshold 0.95: repeated_count 0 -> 4\n$t.Cell(21, 5).Range.Text = \"False\" # threshold 0.95: no_repeats_bool True -> False\n\n$t.Cell(22, 2).Range.Text = \"10\"   # threshold 1.0: row_count 18 -> 10\n$t.Cell(22, 3).Range.Text = \"7\"    # threshold 1.0: unique_count 17 -> 7\n$t.Cell(22, 4).Range.Text = \"5\"    # threshold 1.0: repeated_count 2 -> 5\n\n# Update the two \"Chosen thresholds\" paragraphs via Find & Replace.\n$find1 = $d.Content.Find\n$find1.Text = \"Full Analysis threshold = 0.95.\"\n$find1.Replacement.Text = \"Full Analysis threshold = 0.90.\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n$find2 = $d.Content.Find\n$find2.Text = \"Core-Level Analysis threshold = 0.65.\"\n$find2.Replacement.Text = \"Core-Level Analysis: NONE found (no threshold yields no repeats & \u22649 rows?).\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
